# Quarterly indexing esoteric bug-fix operation
# Column A (rows 2-73) holds an index date per row. Each date was the
# 1st of a quarter-start month; fix it to be the 15th of the *next*
# calendar month (off-by-one-month indexing bug).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excelEpoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($serial -ne $null) {
        $oldDate = $excelEpoch.AddDays([double]$serial)
        $bumped = $oldDate.AddMonths(1)
        $newDate = Get-Date -Year $bumped.Year -Month $bumped.Month -Day 15 -Hour 0 -Minute 0 -Second 0
        $cell.Value2 = $newDate.ToOADate()
    }
}
